$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2016
$ws.Range("D2").Value = 100

$ws.Range("A3").Value = 2016
$ws.Range("D3").Value = 120

$ws.Range("A4").Value = 2017
$ws.Range("D4").Value = 200

$ws.Range("A5").Value = 2017
$ws.Range("D5").Value = 150

$ws.Range("A6").Value = 2018
$ws.Range("D6").Value = 280

$ws.Range("A7").Value = 2018
$ws.Range("D7").Value = 250
